# Update the daily shortage report: refresh the "current balance" (H column)
# figures for several items and bump the generated-at timestamp footer.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7  - ABIMOL 300MG 5 RECTAL SUPP.          6:0  -> 7:0
$ws.Range("H7").Value = "7:0"

# Row 9  - CETAL 250MG/5ML 60ML SUSP            19:0 -> 20:0
$ws.Range("H9").Value = "20:0"

# Row 15 - OPLEX-N SYRUP 125ML                  4:0  -> 5:0
$ws.Range("H15").Value = "5:0"

# Row 16 - PANTOPI 40MG 14 TAB                  0:0  -> 0:1
$ws.Range("H16").Value = "0:1"

# Row 18 - T4-THYRO 50 MCG 100 TABS.            3:0  -> 4:0
$ws.Range("H18").Value = "4:0"

# Row 19 - TEGRETOL CR 400MG 20 F.C. DIVITABS   0:1  -> 1:0
$ws.Range("H19").Value = "1:0"

# Row 24 - صوفي طويل جدا                         15:0 -> 16:0
$ws.Range("H24").Value = "16:0"

# Footer timestamp (A27): bump generation time by a minute
$ws.Range("A27").Value = "Thursday, 29 May, 2025 1:05 PM"
